# Daily attendance processing - 2025-12-19 08:37:14
#
# Normalises the "Recorded By" column (column G) on the active sheet:
# whenever the comma-separated list of recorders begins with the literal
# token "System, ", that leading "System" entry is moved from the front of
# the list to the end (e.g. "System, foo@bar.com" -> "foo@bar.com, System").
# Entries that do not start with "System, " (including bare "System" with
# no other recorder, or lists where "System"/"system" already appear later)
# are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the "Recorded By" column from the header row (row 1) instead of
# hard-coding a column index, and work over the sheet's used range so the
# script keeps working even if rows/columns are rearranged.
$used = $ws.UsedRange
$firstRow = $used.Row
$firstCol = $used.Column
$lastRow = $firstRow + $used.Rows.Count - 1
$lastCol = $firstCol + $used.Columns.Count - 1

$recordedByCol = 0
for ($c = $firstCol; $c -le $lastCol; $c++) {
    $header = $ws.Cells.Item($firstRow, $c).Text
    if ($header -eq "Recorded By") {
        $recordedByCol = $c
        break
    }
}

if ($recordedByCol -eq 0) {
    # Fallback to the known column (G) if the header could not be found.
    $recordedByCol = 7
}

$prefix = "System, "
$changed = 0

for ($r = $firstRow + 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $recordedByCol)
    $val = $cell.Text

    if ([string]::IsNullOrEmpty($val)) {
        continue
    }

    if ($val.StartsWith($prefix)) {
        $rest = $val.Substring($prefix.Length)
        $newVal = $rest + ", System"
        $cell.Value = $newVal
        $changed++
    }
}

Write-Output "Reordered 'System' to the end in $changed 'Recorded By' cell(s)."
